$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1895043731778426
$ws.Range("C2").Value = 0.5422740524781341
$ws.Range("J2").Value = 0.01457725947521866
$ws.Range("P2").Value = 0.1516034985422741
$ws.Range("S2").Value = 0.1020408163265306
$ws.Range("B3").Value = 0.01507537688442211
$ws.Range("C3").Value = 0.03015075376884422
$ws.Range("J3").Value = 0.02512562814070352
$ws.Range("P3").Value = 0.7135678391959799
$ws.Range("S3").Value = 0.2160804020100502
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.6956521739130435
$ws.Range("S4").Value = 0.2608695652173913
$ws.Range("B6").Value = 0.07423580786026202
$ws.Range("D6").Value = 0.008733624454148471
$ws.Range("F6").Value = 0.07860262008733625
$ws.Range("J6").Value = 0.2532751091703057
$ws.Range("O6").Value = 0.02183406113537118
$ws.Range("Q6").Value = 0.1179039301310044
$ws.Range("R6").Value = 0.07423580786026202
$ws.Range("S6").Value = 0.37117903930131
$ws.Range("B7").Value = 0.1031390134529148
$ws.Range("D7").Value = 0.01345291479820628
$ws.Range("F7").Value = 0.05381165919282511
$ws.Range("J7").Value = 0.09417040358744394
$ws.Range("O7").Value = 0.0179372197309417
$ws.Range("Q7").Value = 0.1748878923766816
$ws.Range("R7").Value = 0.08520179372197309
$ws.Range("S7").Value = 0.4573991031390134
$ws.Range("B8").Value = 0.07481751824817519
$ws.Range("D8").Value = 0.0218978102189781
$ws.Range("F8").Value = 0.06204379562043796
$ws.Range("J8").Value = 0.1222627737226277
$ws.Range("O8").Value = 0.01094890510948905
$ws.Range("Q8").Value = 0.1678832116788321
$ws.Range("R8").Value = 0.1003649635036496
$ws.Range("S8").Value = 0.4397810218978102
$ws.Range("B9").Value = 0.0963855421686747
$ws.Range("D9").Value = 0.01807228915662651
$ws.Range("F9").Value = 0.03614457831325301
$ws.Range("J9").Value = 0.07228915662650602
$ws.Range("O9").Value = 0.03614457831325301
$ws.Range("Q9").Value = 0.2228915662650602
$ws.Range("R9").Value = 0.108433734939759
$ws.Range("S9").Value = 0.4096385542168675
$ws.Range("B10").Value = 0.1234567901234568
$ws.Range("D10").Value = 0.0196078431372549
$ws.Range("E10").Value = 0.0007262164124909223
$ws.Range("F10").Value = 0.06899055918663761
$ws.Range("J10").Value = 0.1270878721859114
$ws.Range("O10").Value = 0.01888162672476398
$ws.Range("Q10").Value = 0.1895424836601307
$ws.Range("R10").Value = 0.08496732026143791
$ws.Range("S10").Value = 0.3667392883079157
$ws.Range("G11").Value = 0.134375
$ws.Range("J11").Value = 0.08125
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.575
$ws.Range("S11").Value = 0.009375
$ws.Range("G12").Value = 0.7842105263157895
$ws.Range("J12").Value = 0.1473684210526316
$ws.Range("K12").Value = 0.01578947368421053
$ws.Range("L12").Value = 0.02631578947368421
$ws.Range("S12").Value = 0.02631578947368421
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.04032258064516129
$ws.Range("H15").Value = 0.1451612903225807
$ws.Range("I15").Value = 0.06048387096774194
$ws.Range("J15").Value = 0.3588709677419355
$ws.Range("K15").Value = 0.06048387096774194
$ws.Range("M15").Value = 0.01209677419354839
$ws.Range("N15").Value = 0.004032258064516129
$ws.Range("O15").Value = 0.04838709677419355
$ws.Range("S15").Value = 0.2701612903225806
$ws.Range("F16").Value = 0.01376146788990826
$ws.Range("H16").Value = 0.1788990825688073
$ws.Range("I16").Value = 0.07798165137614679
$ws.Range("J16").Value = 0.4174311926605505
$ws.Range("K16").Value = 0.07339449541284404
$ws.Range("M16").Value = 0.03211009174311927
$ws.Range("O16").Value = 0.05045871559633028
$ws.Range("S16").Value = 0.1559633027522936
$ws.Range("F17").Value = 0.01106194690265487
$ws.Range("H17").Value = 0.2146017699115044
$ws.Range("I17").Value = 0.06415929203539823
$ws.Range("J17").Value = 0.4292035398230089
$ws.Range("K17").Value = 0.08628318584070796
$ws.Range("M17").Value = 0.01548672566371681
$ws.Range("N17").Value = 0.002212389380530973
$ws.Range("O17").Value = 0.05973451327433629
$ws.Range("S17").Value = 0.1172566371681416
$ws.Range("F18").Value = 0.008849557522123894
$ws.Range("H18").Value = 0.1902654867256637
$ws.Range("I18").Value = 0.07079646017699115
$ws.Range("J18").Value = 0.4070796460176991
$ws.Range("K18").Value = 0.1150442477876106
$ws.Range("M18").Value = 0.02212389380530973
$ws.Range("N18").Value = 0.004424778761061947
$ws.Range("O18").Value = 0.08849557522123894
$ws.Range("S18").Value = 0.09292035398230089
$ws.Range("F19").Value = 0.01465457083042568
$ws.Range("H19").Value = 0.2337752965806001
$ws.Range("I19").Value = 0.06350314026517795
$ws.Range("J19").Value = 0.3642707606420097
$ws.Range("K19").Value = 0.1088625261688765
$ws.Range("M19").Value = 0.01744591765526867
$ws.Range("N19").Value = 0.0006978367062107466
$ws.Range("O19").Value = 0.06908583391486392
$ws.Range("S19").Value = 0.1277041172365666
